# Disaggregation of commodity Copper
# Swap the "Photovoltaic plants" / "Onshore wind plants" rows (row 5 and row 6)
# on every year sheet: the labels in C5/C6 swap, and the EU27+UK values in
# E5/E6 swap accordingly (columns D, F, G are always 0 for these rows).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $label5 = $ws.Range("C5").Value2
    $label6 = $ws.Range("C6").Value2
    $val5 = $ws.Range("E5").Value2
    $val6 = $ws.Range("E6").Value2

    $ws.Range("C5").Value2 = $label6
    $ws.Range("C6").Value2 = $label5
    $ws.Range("E5").Value2 = $val6
    $ws.Range("E6").Value2 = $val5
}
